$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column R values for year 2021
$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 18.953297329007047

# Copy style from Q4/Q5 neighbours as appropriate
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("H5").Copy()
$ws.Range("R5").PasteSpecial(-4122)

# Update selection to match target (Q8 instead of Q9)
$ws.Range("Q8").Select()
